# Apply the "Add files via upload" edit to the Acceptance Test Plan workbook.
# Target sheet is "Test Plan" (the second worksheet); "Instructions" is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")

# --- Rows 2-8: mark Pass, add "EY; 7/8" comment ---
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 5).Value = "Pass"
    $ws.Cells.Item($r, 6).Value = "EY; 7/8"
}

# --- Rows 9-15: mark Pass, add "TJ; 7/8" comment ---
for ($r = 9; $r -le 15; $r++) {
    $ws.Cells.Item($r, 5).Value = "Pass"
    $ws.Cells.Item($r, 6).Value = "TJ; 7/8"
}

# --- Row 16: mark Pass, add "TJ;7/2" comment ---
$ws.Cells.Item(16, 5).Value = "Pass"
$ws.Cells.Item(16, 6).Value = "TJ;7/2"

# --- New user stories / acceptance criteria rows 17-32 ---

# Row 17: Update Gameboard
$ws.Cells.Item(17, 1).Value = "Update Gameboard"
$ws.Cells.Item(17, 2).Value = "Given I am red player when I enter game, I expect to see that I can move my piece first"
$ws.Cells.Item(17, 3).Value = "Pass"
$ws.Cells.Item(17, 5).Value = "Pass"
$ws.Cells.Item(17, 6).Value = "TJ;7/2"

# Row 18
$ws.Cells.Item(18, 2).Value = "Given it's my turn when I drag a piece, I expect to place my piece on any free black space"
$ws.Cells.Item(18, 5).Value = "Pass"
$ws.Cells.Item(18, 6).Value = "TJ;7/2"

# Row 19
$ws.Cells.Item(19, 2).Value = "Given it's my turn when I drag a piece to a spot, I expect to have the option to backup or submit"
$ws.Cells.Item(19, 5).Value = "Fail"
$ws.Cells.Item(19, 6).Value = 'TJ;7/2; Submit Move not working, making subsequent stories all fail. Validate move works though with "fake moves" manually passed in instead of via the board.'

# Row 20
$ws.Cells.Item(20, 2).Value = "Given I have made a valid move when I click submit, I expect it to be the opponents turn"
$ws.Cells.Item(20, 5).Value = "Fail"
$ws.Cells.Item(20, 6).Value = "TJ; 7/6"

# Row 21
$ws.Cells.Item(21, 2).Value = "Given I have made an invalid move when I click submit, I expect an error message to appear stating why the move was invalid"
$ws.Cells.Item(21, 5).Value = "Fail"
$ws.Cells.Item(21, 6).Value = "TJ; 7/6"

# Row 22
$ws.Cells.Item(22, 2).Value = "Given I am waiting for my turn when the opponent submits their move, I expect for my board to refresh to a new board and to be able to move pieces"
$ws.Cells.Item(22, 5).Value = "Fail"
$ws.Cells.Item(22, 6).Value = "TJ; 7/6"

# Row 23: Backup
$ws.Cells.Item(23, 1).Value = "Backup"
$ws.Cells.Item(23, 2).Value = "Given I have made a move when I click backup, I expect my board to display the original orientation. "
$ws.Cells.Item(23, 5).Value = "Fail"
$ws.Cells.Item(23, 6).Value = "TJ; 7/11"

# Row 24: Resignation
$ws.Cells.Item(24, 1).Value = "Resignation"
$ws.Cells.Item(24, 2).Value = "Given it's my turn when I haven't touched any pieces, I expect to be able to resign"
$ws.Cells.Item(24, 5).Value = "Fail"
$ws.Cells.Item(24, 6).Value = "TJ; 7/11"

# Row 25: Leaving the Game
$ws.Cells.Item(25, 1).Value = "Leaving the Game"
$ws.Cells.Item(25, 2).Value = "Given that I have left the game, I expect to be reconnected to the game if it's still in progress"
$ws.Cells.Item(25, 5).Value = "Fail"
$ws.Cells.Item(25, 6).Value = "TJ; 7/11"

# Row 26: Generic Move
$ws.Cells.Item(26, 1).Value = "Generic Move"
$ws.Cells.Item(26, 2).Value = "Given that I made a generic move, I expect the move to be appropriately checked before being made."
$ws.Cells.Item(26, 5).Value = "Fail"
$ws.Cells.Item(26, 6).Value = "TJ; 7/11"

# Row 27: Capturing Pieces
$ws.Cells.Item(27, 1).Value = "Capturing Pieces"
$ws.Cells.Item(27, 2).Value = "Given that when I jump a piece, I expect that piece to be removed from play."
$ws.Cells.Item(27, 5).Value = "Fail"
$ws.Cells.Item(27, 6).Value = "TJ; 7/11"

# Row 28
$ws.Cells.Item(28, 2).Value = "Given that I have a jump, I expect that I will be forced to make a jump move."
$ws.Cells.Item(28, 5).Value = "Fail"
$ws.Cells.Item(28, 6).Value = "TJ; 7/11"

# Row 29: Promoting to King
$ws.Cells.Item(29, 1).Value = "Promoting to King"
$ws.Cells.Item(29, 2).Value = "Given that I move my piece to the other side, I expect to be able to promote into a King."
$ws.Cells.Item(29, 5).Value = "Fail"
$ws.Cells.Item(29, 6).Value = "TJ; 7/11"

# Row 30
$ws.Cells.Item(30, 2).Value = "Given that when I have a King piece, I can move and jump backwards."
$ws.Cells.Item(30, 5).Value = "Fail"
$ws.Cells.Item(30, 6).Value = "TJ; 7/11"

# Row 31: Multiple Games
$ws.Cells.Item(31, 1).Value = "Multiple Games"
$ws.Cells.Item(31, 2).Value = "Given that there are a bunch of players in the lobby, there can be more than one game running at a time."

# Row 32
$ws.Cells.Item(32, 2).Value = "Given that players are done with a game, they are automatically returned to lobby and available to be challenged again."

# --- Remove the two now-unused trailing blank rows (596, 597) ---
$ws.Rows.Item(597).Delete()
$ws.Rows.Item(596).Delete()

# --- Restore view/selection state ---
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("B32").Select()
